$wb = $excel.ActiveWorkbook

# Sheet 1: Linear Regressor
$ws1 = $wb.Worksheets.Item("Linear Regressor")
$ws1.Range("B2").Value = 1115.192198589515
$ws1.Range("C2").Value = 33.39449353695179
$ws1.Range("D2").Value = 4.955288999994352

$ws1.Range("B3").Value = 465.8193664687402
$ws1.Range("C3").Value = 21.58284889602714
$ws1.Range("D3").Value = 2.341955098319418

$ws1.Range("B4").Value = 2454.969604111337
$ws1.Range("C4").Value = 49.5476498343901
$ws1.Range("D4").Value = 1.00335087302598

$ws1.Range("B5").Value = 1827.416830215084
$ws1.Range("C5").Value = 42.74829622587413
$ws1.Range("D5").Value = 1.160871396469

$ws1.Range("B6").Value = 3004.840114938119
$ws1.Range("C6").Value = 54.81642194578299
$ws1.Range("D6").Value = 0.9239516363625029

# Sheet 2: Random Forest
$ws2 = $wb.Worksheets.Item("Random Forest")
$ws2.Range("B2").Value = 1505.099261808911
$ws2.Range("C2").Value = 38.79560879544115
$ws2.Range("D2").Value = 6.687817423197884

$ws2.Range("B3").Value = 1016.529174610843
$ws2.Range("C3").Value = 31.88305466248243
$ws2.Range("D3").Value = 5.110705682156411

$ws2.Range("B4").Value = 1482.690546430488
$ws2.Range("C4").Value = 38.50572095715764
$ws2.Range("D4").Value = 0.6059785227878198

$ws2.Range("B5").Value = 1162.61470187869
$ws2.Range("C5").Value = 34.09713627093469
$ws2.Range("D5").Value = 0.7385540781992583

$ws2.Range("B6").Value = 1983.703941210413
$ws2.Range("C6").Value = 44.53879142063033
$ws2.Range("D6").Value = 0.6099647343725177

# Sheet 3: XGBoost
$ws3 = $wb.Worksheets.Item("XGBoost")
$ws3.Range("B2").Value = 2014.03388078173
$ws3.Range("C2").Value = 44.87798882282639
$ws3.Range("D2").Value = 8.949237582253987

$ws3.Range("B3").Value = 1602.105022537301
$ws3.Range("C3").Value = 40.02630413287368
$ws3.Range("D3").Value = 8.054748891223232

$ws3.Range("B4").Value = 1965.705748943596
$ws3.Range("C4").Value = 44.33628027861151
$ws3.Range("D4").Value = 0.8033877796335012

$ws3.Range("B5").Value = 1737.856715191586
$ws3.Range("C5").Value = 41.68760865283095
$ws3.Range("D5").Value = 1.103978095457313

$ws3.Range("B6").Value = 2073.3746820848
$ws3.Range("C6").Value = 45.53432421904161
$ws3.Range("D6").Value = 0.6375373920167099
